# Update iServ stats for 2026-01 (row 26)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6531
$ws.Range("C26").Value = 1019
$ws.Range("D26").Value = 6084592
$ws.Range("E26").Value = 931.6478334098913
$ws.Range("F26").Value = 10.37688017576475
$ws.Range("G26").Value = 8.174097664543535
$ws.Range("H26").Value = 26.71183338806846
